# Product Backlog with some updates, still needs review
# Applies the Sprint 1 sheet updates: revised "Hours Worked" figures, a new
# "Not needed" note, a couple of highlighted rows, and re-grouped I-column
# formulas (Excel turns a dragged/re-entered formula range into a shared
# formula group automatically).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sprint 1")

# ---------------------------------------------------------------------
# 1. Updated "Hours Worked" (column H) figures
# ---------------------------------------------------------------------
$ws.Range("H22").Value = 5
$ws.Range("H23").Value = 4
$ws.Range("H24").Value = 3

$ws.Range("H52").Value = 5
$ws.Range("H53").Value = 4
$ws.Range("H54").Value = 2

# Row 79 had no "Hours Worked" entry yet - give it an explicit 0.
$ws.Range("H79").Value = 0

# ---------------------------------------------------------------------
# 2. Re-enter the "Hours Left" (column I) formulas per task block so
#    Excel groups each contiguous block into a shared formula, matching
#    how the existing I4:I11 / I72:I74 groups were already built.
# ---------------------------------------------------------------------
$ws.Range("I20:I28").Formula = "=G20-H20"
$ws.Range("I36:I42").Formula = "=G36-H36"
$ws.Range("I44:I49").Formula = "=G44-H44"
$ws.Range("I51:I58").Formula = "=G51-H51"
$ws.Range("I60:I66").Formula = "=G60-H60"

# ---------------------------------------------------------------------
# 3. Flag a cut task ("Simple UI" sub-item) as not needed, in red.
# ---------------------------------------------------------------------
$ws.Range("E21:I21").Font.Color = 255          # RGB(255,0,0) -> red
$ws.Range("J21").Value = "Not needed"

# ---------------------------------------------------------------------
# 4. Highlight the over-budget task row in green.
# ---------------------------------------------------------------------
$ws.Range("E46:I46").Font.Color = 5287936      # RGB(0,176,80) -> green

# ---------------------------------------------------------------------
# 5. Move the frozen-pane scroll position / active selection down to the
#    rows now being worked on.
# ---------------------------------------------------------------------
$ws.Activate()
$ws.Range("A54").Select()
$excel.ActiveWindow.FreezePanes = $true
$ws.Range("J60").Select()
